$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '61.150.93'
$ws.Cells.Item(2, 5).Value = '  +0.85%  '

$ws.Cells.Item(3, 4).Value = '2.659.42'
$ws.Cells.Item(3, 5).Value = '  +1.59%  '

$ws.Cells.Item(4, 5).Value = '  +0.00%  '

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '533.78'
$ws.Cells.Item(5, 5).Value = '  +4.24%  '

$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '156.75'
$ws.Cells.Item(6, 5).Value = '  +1.03%  '

$ws.Cells.Item(7, 5).Value = '  +0.02%  '

$ws.Cells.Item(8, 5).Value = '  +0.90%  '

$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '6.59'
$ws.Cells.Item(9, 5).Value = '  -1.62%  '

$ws.Cells.Item(10, 5).Value = '  +5.06%  '

$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.353'
$ws.Cells.Item(11, 5).Value = '  +1.72%  '

$ws.Cells.Item(12, 5).Value = '  -0.03%  '

$ws.Cells.Item(13, 4).Value = '3.121.94'
$ws.Cells.Item(13, 5).Value = '  +1.52%  '

$ws.Cells.Item(14, 4).Value = '61.127.25'
$ws.Cells.Item(14, 5).Value = '  +0.94%  '

$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '22.11'
$ws.Cells.Item(15, 5).Value = '  +2.08%  '

$ws.Cells.Item(16, 5).Value = '  +2.40%  '

$ws.Cells.Item(17, 4).Value = '2.658.04'
$ws.Cells.Item(17, 5).Value = '  +1.20%  '

$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '355.36'
$ws.Cells.Item(19, 5).Value = '  +1.16%  '

$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '10.72'
$ws.Cells.Item(20, 5).Value = '  +0.89%  '

$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '6.28'
$ws.Cells.Item(21, 5).Value = '  +1.57%  '

$ws.Cells.Item(22, 5).Value = '  +0.55%  '

$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '61.65'
$ws.Cells.Item(23, 5).Value = '  +1.70%  '

$ws.Cells.Item(24, 5).Value = '  +2.08%  '

$ws.Cells.Item(25, 5).Value = '  +1.20%  '

$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '1.00'
$ws.Cells.Item(26, 5).Value = '  +0.63%  '

$ws.Cells.Item(27, 4).Value = '0.0₃0864'
$ws.Cells.Item(27, 5).Value = '  +2.11%  '

$ws.Cells.Item(28, 5).Value = '  +0.00%  '

$ws.Cells.Item(29, 5).Value = '  +0.00%  '

$ws.Cells.Item(31, 5).Value = '  +3.96%  '

$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '150.25'
$ws.Cells.Item(33, 5).Value = '  -0.11%  '

$ws.Cells.Item(34, 5).Value = '  +3.70%  '

$ws.Cells.Item(35, 5).Value = '  +0.88%  '

$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '0.925'
$ws.Cells.Item(36, 5).Value = '  +8.99%  '

$ws.Cells.Item(37, 5).Value = '  -0.48%  '

$ws.Cells.Item(38, 2).Value = 'Stacks'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '1.50'
$ws.Cells.Item(38, 5).Value = '  +0.17%  '

$ws.Cells.Item(39, 2).Value = 'Filecoin'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '3.83'
$ws.Cells.Item(39, 5).Value = '  +1.09%  '

$ws.Cells.Item(40, 2).Value = 'Bittensor'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '306.13'
$ws.Cells.Item(40, 5).Value = '  +4.01%  '

$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.653'
$ws.Cells.Item(41, 5).Value = '  +4.07%  '

$ws.Cells.Item(42, 5).Value = '  +1.86%  '

$ws.Cells.Item(43, 5).Value = '  +1.97%  '

$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '20.20'
$ws.Cells.Item(44, 5).Value = '  +1.36%  '

$ws.Cells.Item(45, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '0.997'
$ws.Cells.Item(45, 5).Value = '  +0.07%  '

$ws.Cells.Item(46, 2).Value = 'RenderToken'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '5.05'
$ws.Cells.Item(46, 5).Value = '  +2.77%  '

$ws.Cells.Item(47, 5).Value = '  +2.36%  '

$ws.Cells.Item(48, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '19.18'
$ws.Cells.Item(48, 5).Value = '  +7.67%  '

$ws.Cells.Item(49, 2).Value = 'WhiteBITCoin'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '10.35'
$ws.Cells.Item(49, 5).Value = '  +0.26%  '

$ws.Cells.Item(50, 4).Value = '1.993.98'
$ws.Cells.Item(50, 5).Value = '  -0.65%  '

$ws.Cells.Item(51, 5).Value = '  +1.70%  '
